# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 12.0302756157461)
    3  = @(0.3464964993005633, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 1.611428759079651)
    4  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    5  = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 1.705647867635037)
    6  = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 0, 22.31973251085698)
    7  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    8  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    9  = @(0.3464964993005633, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 1.611428759079651)
    10 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
